# "complete monthly and re-run daily"
# The underlying raw counts for a handful of records increased by 1 (an
# extra event was counted for the re-run), which ripples through the
# per-record ratio (S/T), the per-day-per-person rollup (F/G) and the two
# productivity summary sheets (D) that are derived from it.

$wb = $excel.ActiveWorkbook

# --- Sheet "team_df": per-shift-record detail --------------------------
$ws1 = $wb.Worksheets.Item("team_df")

# Row 17 (sp99018 / 張宜君)
$ws1.Range("S17").Value = 7
$ws1.Range("T17").Value = 46
$ws1.Range("U17").Value = 0.1521739130434783

# Row 18 (sp99063 / 林川評)
$ws1.Range("S18").Value = 7
$ws1.Range("T18").Value = 46
$ws1.Range("U18").Value = 0.1521739130434783

# Row 23 (sp99004 / 陳衣玲)
$ws1.Range("S23").Value = 6
$ws1.Range("T23").Value = 40
$ws1.Range("U23").Value = 0.15

# --- Sheet "team_df_day": per-person-per-day rollup ---------------------
$ws2 = $wb.Worksheets.Item("team_df_day")

# Row 2 (sp99004 / 陳衣玲)
$ws2.Range("F2").Value = 14
$ws2.Range("G2").Value = 107
$ws2.Range("H2").Value = 0.1308411214953271

# Row 5 (sp99018 / 張宜君)
$ws2.Range("F5").Value = 12
$ws2.Range("G5").Value = 95
$ws2.Range("H5").Value = 0.1263157894736842

# Row 10 (sp99063 / 林川評)
$ws2.Range("F10").Value = 12
$ws2.Range("G10").Value = 95
$ws2.Range("H10").Value = 0.1263157894736842

# --- Sheet "productivity_tl": TL productivity summary -------------------
$ws3 = $wb.Worksheets.Item("productivity_tl")

$ws3.Range("D2").Value = 0.1308411214953271
$ws3.Range("D5").Value = 0.1263157894736842
$ws3.Range("D10").Value = 0.1263157894736842

# --- Sheet "productivity_team_function": team/function summary ----------
$ws4 = $wb.Worksheets.Item("productivity_team_function")

$ws4.Range("D2").Value = 0.1308411214953271
$ws4.Range("D5").Value = 0.1263157894736842
$ws4.Range("D10").Value = 0.1263157894736842

$wb.Save()
